$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. '27.160.34', '306.55'); force
# the whole Price column to Text format first so COM doesn't coerce these
# assignments into floating-point numbers, then restore the default style.
$ws.Range("D2:D51").NumberFormat = "@"

# Price (D) updates
$ws.Range('D2').Value = '27.160.34'
$ws.Range('D3').Value = '1.901.16'
$ws.Range('D5').Value = '306.55'
$ws.Range('D7').Value = '0.5233'
$ws.Range('D8').Value = '0.3765'
$ws.Range('D9').Value = '0.07243'
$ws.Range('D11').Value = '0.8981'
$ws.Range('D12').Value = '0.08378'
$ws.Range('D13').Value = '1.909.77'
$ws.Range('D14').Value = '94.56'
$ws.Range('D15').Value = '5.259'
$ws.Range('D17').Value = '0.000008581'
$ws.Range('D18').Value = '14.48'
$ws.Range('D20').Value = '27.197.58'
$ws.Range('D22').Value = '2.149.37'
$ws.Range('D24').Value = '6.413'
$ws.Range('D25').Value = '146.55'
$ws.Range('D26').Value = '2.278'
$ws.Range('D28').Value = '18.12'
$ws.Range('D29').Value = '114.61'
$ws.Range('D30').Value = '4.916'
$ws.Range('D31').Value = '4.780'
$ws.Range('D32').Value = '0.09224'
$ws.Range('D33').Value = '0.8187'
$ws.Range('D34').Value = '0.05049'
$ws.Range('D35').Value = '1.233'
$ws.Range('D36').Value = '2.963'
$ws.Range('D37').Value = '3.356'
$ws.Range('D38').Value = '2.564'
$ws.Range('D39').Value = '0.5680'
$ws.Range('D40').Value = '0.01971'
$ws.Range('D42').Value = '6.648'
$ws.Range('D43').Value = '8.908'
$ws.Range('D44').Value = '118.10'
$ws.Range('D46').Value = '0.4815'
$ws.Range('D47').Value = '1.001'
$ws.Range('D48').Value = '10.13'
$ws.Range('D49').Value = '1.608'
$ws.Range('D51').Value = '63.54'

$ws.Range("D2:D51").Style = "Normal"

# Coin / Link / Volume(1h) updates
$ws.Range('E2').Value = '  +0.84%  '
$ws.Range('E3').Value = '  +1.39%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('E5').Value = '  +0.06%  '
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('E7').Value = '  +1.35%  '
$ws.Range('E8').Value = '  +1.41%  '
$ws.Range('E9').Value = '  +0.74%  '
$ws.Range('E10').Value = '  +2.24%  '
$ws.Range('E11').Value = '  +0.00%  '
$ws.Range('E12').Value = '  +10.85%  '
$ws.Range('E13').Value = '  +2.36%  '
$ws.Range('E14').Value = '  -0.33%  '
$ws.Range('E15').Value = '  +0.16%  '
$ws.Range('E16').Value = '  +0.09%  '
$ws.Range('E17').Value = '  +1.13%  '
$ws.Range('E18').Value = '  +1.74%  '
$ws.Range('E19').Value = '  +0.12%  '
$ws.Range('E20').Value = '  +0.88%  '
$ws.Range('E22').Value = '  +2.63%  '
$ws.Range('E23').Value = '  +1.71%  '
$ws.Range('E24').Value = '  -0.20%  '
$ws.Range('E25').Value = '  +0.42%  '
$ws.Range('E26').Value = '  +7.80%  '
$ws.Range('E27').Value = '  -1.60%  '
$ws.Range('E28').Value = '  +0.52%  '
$ws.Range('E30').Value = '  +0.28%  '
$ws.Range('E31').Value = '  +0.86%  '
$ws.Range('E32').Value = '  +0.54%  '
$ws.Range('E33').Value = '  +8.77%  '
$ws.Range('E34').Value = '  +0.39%  '
$ws.Range('E35').Value = '  +5.37%  '
$ws.Range('E36').Value = '  -0.81%  '
$ws.Range('E37').Value = '  +2.34%  '
$ws.Range('E38').Value = '  +3.08%  '
$ws.Range('E39').Value = '  +2.09%  '
$ws.Range('E40').Value = '  -0.96%  '
$ws.Range('E41').Value = '  +0.27%  '
$ws.Range('E42').Value = '  +1.14%  '
$ws.Range('E43').Value = '  +2.17%  '
$ws.Range('E44').Value = '  +1.87%  '
$ws.Range('E45').Value = '  +0.68%  '
$ws.Range('E46').Value = '  +1.14%  '
$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('E47').Value = '  +0.17%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('E48').Value = '  +0.73%  '
$ws.Range('E49').Value = '  +3.02%  '
$ws.Range('E50').Value = '  +0.90%  '
$ws.Range('E51').Value = '  +0.37%  '
